{"js": "// Replace the three-digit-number-divided-by-one-digit-number answers\n// throughout the document's table cells with the new generated values.\n// Each \"find\" string is unique in the document, so searching the whole\n// body and replacing the single match is unambiguous per pair.\n\nconst replacements = [\n  { find: \"757\u00f79=84, 1\", replace: \"302\u00f76=50, 2\" },\n  { find: \"758\u00f79=84, 2\", replace: \"463\u00f72=231, 1\" },\n  { find: \"221\u00f75=44, 1\", replace: \"977\u00f73=325, 2\" },\n  { find: \"608\u00f75=121, 3\", replace: \"362\u00f75=72, 2\" },\n  { find: \"475\u00f76=79, 1\", replace: \"389\u00f77=55, 4\" },\n  { find: \"359\u00f76=59, 5\", replace: \"999\u00f78=124, 7\" },\n  { find: \"831\u00f76=138, 3\", replace: \"203\u00f72=101, 1\" },\n  { find: \"775\u00f79=86, 1\", replace: \"114\u00f75=22, 4\" },\n  { find: \"176\u00f76=29, 2\", replace: \"801\u00f78=100, 1\" },\n  { find: \"599\u00f76=99, 5\", replace: \"332\u00f76=55, 2\" },\n  { find: \"290\u00f72=145, 0\", replace: \"814\u00f79=90, 4\" },\n  { find: \"586\u00f79=65, 1\", replace: \"795\u00f76=132, 3\" },\n  { find: \"308\u00f77=44, 0\", replace: \"522\u00f79=58, 0\" },\n  { find: \"873\u00f77=124, 5\", replace: \"227\u00f76=37, 5\" },\n  { find: \"847\u00f73=282, 1\", replace: \"316\u00f75=63, 1\" },\n  { find: \"968\u00f75=193, 3\", replace: \"842\u00f78=105, 2\" },\n  { find: \"642\u00f76=107, 0\", replace: \"144\u00f77=20, 4\" },\n  { find: \"191\u00f77=27, 2\", replace: \"592\u00f73=197, 1\" },\n  { find: \"407\u00f78=50, 7\", replace: \"369\u00f73=123, 0\" },\n  { find: \"183\u00f74=45, 3\", replace: \"226\u00f74=56, 2\" },\n  { find: \"171\u00f76=28, 3\", replace: \"492\u00f78=61, 4\" },\n  { find: \"154\u00f77=22, 0\", replace: \"120\u00f74=30, 0\" },\n  { find: \"827\u00f72=413, 1\", replace: \"819\u00f77=117, 0\" },\n  { find: \"873\u00f73=291, 0\", replace: \"775\u00f76=129, 1\" },\n  { find: \"227\u00f79=25, 2\", replace: \"137\u00f74=34, 1\" },\n];\n\nconst body = context.document.body;\n\nfor (const { find, replace } of replacements) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replace, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the three-digit-number-divided-by-one-digit-number answers\n# throughout the document's table cells with the new generated values.\n# Each \"find\" string is unique in the document, so a simple\n# Find/Replace (wdReplaceAll) per pair is sufficient and unambiguous.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Find = \"757\u00f79=84, 1\";  Replace = \"302\u00f76=50, 2\" },\n    @{ Find = \"758\u00f79=84, 2\";  Replace = \"463\u00f72=231, 1\" },\n    @{ Find = \"221\u00f75=44, 1\";  Replace = \"977\u00f73=325, 2\" },\n    @{ Find = \"608\u00f75=121, 3\"; Replace = \"362\u00f75=72, 2\" },\n    @{ Find = \"475\u00f76=79, 1\";  Replace = \"389\u00f77=55, 4\" },\n    @{ Find = \"359\u00f76=59, 5\";  Replace = \"999\u00f78=124, 7\" },\n    @{ Find = \"831\u00f76=138, 3\"; Replace = \"203\u00f72=101, 1\" },\n    @{ Find = \"775\u00f79=86, 1\";  Replace = \"114\u00f75=22, 4\" },\n    @{ Find = \"176\u00f76=29, 2\";  Replace = \"801\u00f78=100, 1\" },\n    @{ Find = \"599\u00f76=99, 5\";  Replace = \"332\u00f76=55, 2\" },\n    @{ Find = \"290\u00f72=145, 0\"; Replace = \"814\u00f79=90, 4\" },\n    @{ Find = \"586\u00f79=65, 1\";  Replace = \"795\u00f76=132, 3\" },\n    @{ Find = \"308\u00f77=44, 0\";  Replace = \"522\u00f79=58, 0\" },\n    @{ Find = \"873\u00f77=124, 5\"; Replace = \"227\u00f76=37, 5\" },\n    @{ Find = \"847\u00f73=282, 1\"; Replace = \"316\u00f75=63, 1\" },\n    @{ Find = \"968\u00f75=193, 3\"; Replace = \"842\u00f78=105, 2\" },\n    @{ Find = \"642\u00f76=107, 0\"; Replace = \"144\u00f77=20, 4\" },\n    @{ Find = \"191\u00f77=27, 2\";  Replace = \"592\u00f73=197, 1\" },\n    @{ Find = \"407\u00f78=50, 7\";  Replace = \"369\u00f73=123, 0\" },\n    @{ Find = \"183\u00f74=45, 3\";  Replace = \"226\u00f74=56, 2\" },\n    @{ Find = \"171\u00f76=28, 3\";  Replace = \"492\u00f78=61, 4\" },\n    @{ Find = \"154\u00f77=22, 0\";  Replace = \"120\u00f74=30, 0\" },\n    @{ Find = \"827\u00f72=413, 1\"; Replace = \"819\u00f77=117, 0\" },\n    @{ Find = \"873\u00f73=291, 0\"; Replace = \"775\u00f76=129, 1\" },\n    @{ Find = \"227\u00f79=25, 2\";  Replace = \"137\u00f74=34, 1\" }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.Find\n    $find.Replacement.Text = $r.Replace\n    $find.Execute(\n        $r.Find,    # FindText\n        $false,     # MatchCase\n        $false,     # MatchWholeWord\n        $false,     # MatchWildcards\n        $false,     # MatchSoundsLike\n        $false,     # MatchAllWordForms\n        $true,      # Forward\n        1,          # Wrap (wdFindContinue)\n        $false,     # Format\n        $r.Replace, # ReplaceWith\n        2           # Replace (wdReplaceAll)\n    )\n}\n"}
